# Update the "Test_Yearly" sheet so the Facility Code / Facility Name
# question rows are relabelled as Entity Code / Entity Name, to reflect
# that survey responses are now re-associated with the generic "entity"
# rather than a facility specifically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Yearly")

$ws.Range("D2").Value = "Entity Code"
$ws.Range("D3").Value = "Entity Name"
